# Commit: "change S+LFM to S/LFM in taxonomy"
#
# Fix taxonomy typo "S+LFM/CDL" -> "S/LFM/CDL" (the separate, already-correct
# "S/LFM+CDL" entries must be left untouched) and apply the word-wrap
# formatting that accompanied the fix in the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Correct the taxonomy string everywhere it occurs (Offices/Trade/Hotels
#    summary cells each contain it exactly once).
[void]$ws.Cells.Replace("S+LFM/CDL", "S/LFM/CDL")

# 2. Wrap the long taxonomy text in the data row so it is readable.
$ws.Range("B2:D2").WrapText = $true

# 3. Widen the three data columns to fit the wrapped text.
$ws.Columns.Item(2).ColumnWidth = 41.25
$ws.Columns.Item(3).ColumnWidth = 40.75
$ws.Columns.Item(4).ColumnWidth = 44.6

# 4. Set the row height to accommodate the wrapped, multi-line text.
$ws.Rows.Item(2).RowHeight = 365

# 5. Restore the selection to where the editor left it.
[void]$ws.Range("C5").Select()
